# attendanceTracker/data.xlsx: record the attendance for COE 354.
# RollNo 1 (kobinaholison2002@gmail.com) now shows 2 recorded attendances
# for the "COE 354" column (was 0) - "application is working now".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C2").Value = 2
